# This sheet contains a set of UML-like "connector" cells (arrows) spread across
# columns A-G per row, alongside a label in column H. The edit shifts most of
# these arrow/connector cells from column A over to the right (mostly into
# B/C/E depending on the row), matching the new indentation/nesting levels of
# the corresponding class/method labels added for the research paper source
# listing. Using Range.Cut(destination) both moves the cell's value AND its
# style/formatting in one step, and leaves the source cell blank (but keeps
# its original style) -- exactly matching how this workbook already
# represents "blank" blank-but-styled cells elsewhere (e.g. A19/A20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Cut($ws.Range("E2"))
$ws.Range("A3").Cut($ws.Range("E3"))
$ws.Range("A4").Cut($ws.Range("E4"))
$ws.Range("A5").Cut($ws.Range("E5"))
$ws.Range("A6").Cut($ws.Range("E6"))
$ws.Range("A7").Cut($ws.Range("E7"))
$ws.Range("A8").Cut($ws.Range("E8"))
$ws.Range("B8").Cut($ws.Range("C8"))
$ws.Range("A9").Cut($ws.Range("E9"))
$ws.Range("B9").Cut($ws.Range("C9"))
$ws.Range("A10").Cut($ws.Range("E10"))
$ws.Range("B10").Cut($ws.Range("C10"))
$ws.Range("A11").Cut($ws.Range("E11"))
$ws.Range("A12").Cut($ws.Range("E12"))
$ws.Range("C12").Cut($ws.Range("B12"))
$ws.Range("A13").Cut($ws.Range("E13"))
$ws.Range("C13").Cut($ws.Range("B13"))
$ws.Range("A14").Cut($ws.Range("E14"))
$ws.Range("C14").Cut($ws.Range("B14"))
$ws.Range("A15").Cut($ws.Range("E15"))
$ws.Range("C15").Cut($ws.Range("B15"))
$ws.Range("A16").Cut($ws.Range("E16"))
$ws.Range("C16").Cut($ws.Range("B16"))
$ws.Range("A17").Cut($ws.Range("E17"))
$ws.Range("C17").Cut($ws.Range("B17"))
$ws.Range("E18").Cut($ws.Range("A18"))
$ws.Range("C18").Cut($ws.Range("B18"))
$ws.Range("C19").Cut($ws.Range("B19"))
$ws.Range("C20").Cut($ws.Range("B20"))
$ws.Range("C21").Cut($ws.Range("B21"))
